# "new 'christmas tree' project" - update input parameters for the resistance calculator
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

$ws.Range("B3").Value = 5      # Vcc - napiecie zasilania
$ws.Range("B4").Value = 2      # Ilosc diod rownoleglych
$ws.Range("B5").Value = 2.5    # Vled - napiecie diody/diod szeregowych
$ws.Range("B6").Value = 40     # Ic - prad wszystkich diod (suma)
$ws.Range("B7").Value = 100    # Hfe - wzmocnienie tranzystora

$ws.Range("B7").Select()
